$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.648.67"
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").Value = "2.015.89"
$ws.Range("E3").Value = "  -4.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("E5").Value = "  -3.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5049"
$ws.Range("E7").Value = "  -3.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4265"
$ws.Range("E8").Value = "  -4.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.24"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09253"
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.130"
$ws.Range("E11").Value = "  -3.61%  "
$ws.Range("E12").Value = "  -5.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.133"
$ws.Range("E13").Value = "  -6.90%  "
$ws.Range("D14").Value = "1.990.53"
$ws.Range("E14").Value = "  -5.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.561"
$ws.Range("E15").Value = "  -5.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.02"
$ws.Range("E16").Value = "  -5.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.012"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001129"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06673"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("E20").Value = "  -6.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.008"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.011"
$ws.Range("E22").Value = "  -4.73%  "
$ws.Range("D23").Value = "29.685.27"
$ws.Range("E23").Value = "  -2.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.06"
$ws.Range("E24").Value = "  -4.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.282"
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.81"
$ws.Range("E27").Value = "  -5.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.466"
$ws.Range("E28").Value = "  -5.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.352"
$ws.Range("E29").Value = "  -7.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.94"
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.063"
$ws.Range("E31").Value = "  -7.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.596"
$ws.Range("E32").Value = "  -9.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09994"
$ws.Range("E33").Value = "  -5.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.894"
$ws.Range("E34").Value = "  -5.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.807"
$ws.Range("E35").Value = "  -3.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.646"
$ws.Range("E36").Value = "  -8.27%  "
$ws.Range("E37").Value = "  -6.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.309"
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06409"
$ws.Range("E39").Value = "  -5.77%  "
$ws.Range("E40").Value = "  -6.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.84"
$ws.Range("E41").Value = "  -5.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2085"
$ws.Range("E42").Value = "  -6.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.010"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6384"
$ws.Range("E44").Value = "  -6.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.56"
$ws.Range("E45").Value = "  -7.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.223"
$ws.Range("E46").Value = "  -5.99%  "
$ws.Range("E47").Value = "  -4.63%  "
$ws.Range("E48").Value = "  -3.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07038"
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.139"
$ws.Range("E50").Value = "  -5.11%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.143"
$ws.Range("E51").Value = "  -6.53%  "
